# Add a new "txt_debug" widget row to the Translation sheet.
#
# Commit: "Add new widget: txt_debug, function unimplemented"
#         "Refactor: Change names of some widgets."
#
# The Translation sheet (B:F columns = Widget Id, Typography, Alignment,
# Direction, Text) gets one new row appended right after the last existing
# data row (row 22 -> new row 23), describing the new "res_txt_debug"
# widget whose displayed text is "Debug".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New widget entry, appended as row 23 (directly below the last used row).
$ws.Range("B23").Value = "res_txt_debug"
$ws.Range("C23").Value = "Default"
$ws.Range("D23").Value = "Left"
$ws.Range("E23").Value = "LTR"
$ws.Range("F23").Value = "Debug"

# Touch row 24 (without changing anything visible) so a trailing blank row
# marker is kept below the new data, matching the generated sheet layout.
$ws.Rows.Item(24).OutlineLevel = 0
